$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Update status text "Handed back: in sync with en-US" -> "Ready for handoff"
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("C2").Value = "Ready for handoff"

# Update timestamps
$wsOverview.Range("G2").Value = "2016-08-17 10:56:15"
$wsDeDe.Range("H2").Value = "2016-08-17 10:56:15"
$wsZhCn.Range("H2").Value = "2016-08-17 10:56:10"

# Adjust column widths (Status columns autofit narrower now)
$wsOverview.Range("E1").ColumnWidth = 17.2159881591797
$wsOverview.Range("F1").ColumnWidth = 17.2159881591797
$wsZhCn.Range("C1").ColumnWidth = 17.2159881591797
$wsDeDe.Range("C1").ColumnWidth = 17.2159881591797
